$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 310
$ws.Range("I6").Value = 180
$ws.Range("J6").Value = 526.6667
$ws.Range("K6").Value = 540
$ws.Range("L6").Value = 1580.0001
$ws.Range("M6").Value = -428
$ws.Range("N6").Value = -1804.0001
$ws.Range("H8").Value = 106.333336
$ws.Range("I8").Value = 106.333336
$ws.Range("K8").Value = 319.000008
$ws.Range("M8").Value = -180.000008
$ws.Range("H33").Value = 35968.855
$ws.Range("I33").Value = 40203.12
$ws.Range("K33").Value = 40203.12
$ws.Range("M33").Value = -39974.12
$ws.Range("H40").Value = 1332.6666
$ws.Range("I40").Value = 1300
$ws.Range("J40").Value = 1398
$ws.Range("K40").Value = 1300
$ws.Range("L40").Value = 1398
$ws.Range("M40").Value = -1125
$ws.Range("N40").Value = -1748
$ws.Range("H52").Value = 1367.8572
$ws.Range("I52").Value = 2000
$ws.Range("J52").Value = 1195.4546
$ws.Range("K52").Value = 6000
$ws.Range("L52").Value = 3586.3638
$ws.Range("M52").Value = -5840
$ws.Range("N52").Value = -3906.3638
$ws.Range("H53").Value = 298
$ws.Range("J53").Value = 620.4
$ws.Range("L53").Value = 620.4
$ws.Range("N53").Value = -1894.4
$ws.Range("M69").Value = -17675.375
$ws.Range("H69").Value = 6183.125
$ws.Range("I69").Value = 6183.125
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 18549.375
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""
$ws.Range("H70").Value = 2187.889
$ws.Range("I70").Value = 2298.2
$ws.Range("J70").Value = 2050
$ws.Range("K70").Value = 6894.599999999999
$ws.Range("L70").Value = 6150
$ws.Range("M70").Value = -6624.599999999999
$ws.Range("N70").Value = -6690
$ws.Range("M72").Value = -51280.125
$ws.Range("H72").Value = 6183.125
$ws.Range("I72").Value = 6183.125
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 55648.125
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""
$ws.Range("H73").Value = 2187.889
$ws.Range("I73").Value = 2298.2
$ws.Range("J73").Value = 2050
$ws.Range("K73").Value = 6894.599999999999
$ws.Range("L73").Value = 6150
$ws.Range("M73").Value = -5958.599999999999
$ws.Range("N73").Value = -8022
$ws.Range("H80").Value = 766.7692
$ws.Range("I80").Value = 596.6
$ws.Range("J80").Value = 873.125
$ws.Range("K80").Value = 1789.8
$ws.Range("L80").Value = 2619.375
$ws.Range("M80").Value = -791.8000000000002
$ws.Range("N80").Value = -4615.375
$ws.Range("H83").Value = 766.7692
$ws.Range("I83").Value = 596.6
$ws.Range("J83").Value = 873.125
$ws.Range("K83").Value = 5369.400000000001
$ws.Range("L83").Value = 7858.125
$ws.Range("M83").Value = -377.4000000000005
$ws.Range("N83").Value = -17842.125
$ws.Range("H98").Value = 2764.3572
$ws.Range("I98").Value = 1639.5
$ws.Range("J98").Value = 5576.5
$ws.Range("K98").Value = 1639.5
$ws.Range("L98").Value = 5576.5
$ws.Range("M98").Value = -141.5
$ws.Range("N98").Value = -8572.5
$ws.Range("H116").Value = 1879.3922
$ws.Range("I116").Value = 1379
$ws.Range("J116").Value = 2655
$ws.Range("K116").Value = 1379
$ws.Range("L116").Value = 2655
$ws.Range("M116").Value = 2063
$ws.Range("N116").Value = -9539
$ws.Range("H122").Value = 2764.3572
$ws.Range("I122").Value = 1639.5
$ws.Range("J122").Value = 5576.5
$ws.Range("K122").Value = 4918.5
$ws.Range("L122").Value = 16729.5
$ws.Range("M122").Value = -2468.5
$ws.Range("N122").Value = -21629.5
$ws.Range("H137").Value = 2818
$ws.Range("I137").Value = 2984.4194
$ws.Range("J137").Value = 2244.7778
$ws.Range("K137").Value = 8953.2582
$ws.Range("L137").Value = 6734.3334
$ws.Range("M137").Value = -6403.2582
$ws.Range("N137").Value = -11834.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2408.0476
$ws.Range("I88").Value = 1841.1111
$ws.Range("J88").Value = 2833.25
$ws.Range("K88").Value = 1841.1111
$ws.Range("L88").Value = 2833.25
$ws.Range("M88").Value = -1435.1111
$ws.Range("N88").Value = -3645.25
$ws.Range("H91").Value = 2408.0476
$ws.Range("I91").Value = 1841.1111
$ws.Range("J91").Value = 2833.25
$ws.Range("K91").Value = 1841.1111
$ws.Range("L91").Value = 2833.25
$ws.Range("M91").Value = -437.1111000000001
$ws.Range("N91").Value = -5641.25
$ws.Range("H135").Value = 29993.111
$ws.Range("J135").Value = 29993.111
$ws.Range("L135").Value = 29993.111
$ws.Range("N135").Value = -40133.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1183.3871
$ws.Range("I94").Value = 878
$ws.Range("J94").Value = 3244.75
$ws.Range("K94").Value = 878
$ws.Range("L94").Value = 3244.75
$ws.Range("M94").Value = -427
$ws.Range("N94").Value = -4146.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24199.611
$ws.Range("I31").Value = 39069.414
$ws.Range("J31").Value = 2638.4
$ws.Range("K31").Value = 39069.414
$ws.Range("L31").Value = 2638.4
$ws.Range("M31").Value = -38774.414
$ws.Range("N31").Value = -3228.4
$ws.Range("H34").Value = 24199.611
$ws.Range("I34").Value = 39069.414
$ws.Range("J34").Value = 2638.4
$ws.Range("K34").Value = 39069.414
$ws.Range("L34").Value = 2638.4
$ws.Range("M34").Value = -38867.414
$ws.Range("N34").Value = -3042.4
$ws.Range("H135").Value = 50654.75
$ws.Range("J135").Value = 50654.75
$ws.Range("L135").Value = 50654.75
$ws.Range("N135").Value = -60794.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1110.2972
$ws.Range("I5").Value = 546.55
$ws.Range("J5").Value = 1773.5294
$ws.Range("K5").Value = 1639.65
$ws.Range("L5").Value = 5320.5882
$ws.Range("M5").Value = -1527.65
$ws.Range("N5").Value = -5544.5882
$ws.Range("H40").Value = 48
$ws.Range("I40").Value = 50.88889
$ws.Range("J40").Value = 39.333332
$ws.Range("K40").Value = 203.55556
$ws.Range("L40").Value = 157.333328
$ws.Range("M40").Value = -134.55556
$ws.Range("N40").Value = -295.333328
$ws.Range("H51").Value = 2906.1904
$ws.Range("I51").Value = 1750
$ws.Range("J51").Value = 3027.8948
$ws.Range("K51").Value = 5250
$ws.Range("L51").Value = 9083.6844
$ws.Range("M51").Value = -4790
$ws.Range("N51").Value = -10003.6844
$ws.Range("H132").Value = 1219.0834
$ws.Range("I132").Value = 817.25
$ws.Range("K132").Value = 7355.25
$ws.Range("M132").Value = -4825.25
$ws.Range("H135").Value = 1110.2972
$ws.Range("I135").Value = 546.55
$ws.Range("J135").Value = 1773.5294
$ws.Range("K135").Value = 4918.95
$ws.Range("L135").Value = 15961.7646
$ws.Range("M135").Value = -2383.95
$ws.Range("N135").Value = -21031.7646

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3391.0557
$ws.Range("I80").Value = 2660
$ws.Range("J80").Value = 3672.2307
$ws.Range("K80").Value = 2660
$ws.Range("L80").Value = 3672.2307
$ws.Range("M80").Value = -1662
$ws.Range("N80").Value = -5668.2307
$ws.Range("H83").Value = 3391.0557
$ws.Range("I83").Value = 2660
$ws.Range("J83").Value = 3672.2307
$ws.Range("K83").Value = 13300
$ws.Range("L83").Value = 18361.1535
$ws.Range("M83").Value = -8308
$ws.Range("N83").Value = -28345.1535
$ws.Range("H122").Value = 1847.2142
$ws.Range("I122").Value = 1532.8182
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4598.4546
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2148.4546
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 667.1
$ws.Range("I22").Value = 610.1429000000001
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 610.1429000000001
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -315.1429000000001
$ws.Range("N22").Value = -1390
$ws.Range("H27").Value = 667.1
$ws.Range("I27").Value = 610.1429000000001
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 610.1429000000001
$ws.Range("L27").Value = 800
$ws.Range("M27").Value = -503.1429000000001
$ws.Range("N27").Value = -1014
$ws.Range("H82").Value = 2577.2273
$ws.Range("I82").Value = 1819.8
$ws.Range("J82").Value = 2800
$ws.Range("K82").Value = 1819.8
$ws.Range("L82").Value = 2800
$ws.Range("M82").Value = -1458.8
$ws.Range("N82").Value = -3522
$ws.Range("H85").Value = 2577.2273
$ws.Range("I85").Value = 1819.8
$ws.Range("J85").Value = 2800
$ws.Range("K85").Value = 1819.8
$ws.Range("L85").Value = 2800
$ws.Range("M85").Value = -571.8
$ws.Range("N85").Value = -5296
$ws.Range("H127").Value = 39982.375
$ws.Range("J127").Value = 39982.375
$ws.Range("L127").Value = 39982.375
$ws.Range("N127").Value = -49902.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 39943.332
$ws.Range("J125").Value = 39943.332
$ws.Range("L125").Value = 39943.332
$ws.Range("N125").Value = -49783.332
